# Applies the "cds Primary diagnosis fixed" change:
# Replaces the broken Neo4j/Cypher query used for the "ParticipantsTab" row
# (cell B2 on the "startup" sheet) with the corrected query, adjusts the
# row height to fit the new (longer) text, and moves the active selection
# to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in ['Acinar cell carcinoma']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

# Update the query text for the Participants row.
$ws.Range("B2").Value = $newQuery

# The longer query text needs a taller row to display (matches saved file).
$ws.Rows("2:2").RowHeight = 330.75

# Move the selected/active cell to B2 (was B8).
$ws.Range("B2").Select()
